$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "set.seed(20201210)"
$ws.Range("B4").Value = 3
$ws.Range("C4").Formula = "=68 * B4"
$ws.Range("D4").Value = 399
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 25
$ws.Range("G4").Value = 0.023912037037037034
$ws.Range("H4").Value = 1527.413

$ws.Range("B4").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("C4").NumberFormat = $ws.Range("C2").NumberFormat
$ws.Range("G4").NumberFormat = $ws.Range("G2").NumberFormat

$ws.Range("H4").Select()
